# Cotações atualizadas - 2025-10-25
# Appends the new daily quote row (row 51) to Sheet1, mirroring the
# layout/formatting of the preceding rows (row 50 in particular).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 51
$prevRow = $newRow - 1

# Date serial for 2025-10-25, stored like the other "Data" column cells
# (numeric value, formatted via the same custom date/time number format).
$ws.Cells.Item($newRow, 1).Value = 45955
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($prevRow, 1).NumberFormat

# Quote values are kept as text (comma decimal separator), same as all
# other rows in the sheet.
$ws.Cells.Item($newRow, 2).Value = "21,7048"
$ws.Cells.Item($newRow, 3).Value = "15,5758"
$ws.Cells.Item($newRow, 4).Value = "15,5156"
$ws.Cells.Item($newRow, 5).Value = "15,5156"
